# Update forecast_summary_B08LGKGBKT.xlsx workbook with the corrected
# forecast output:
#   - "Forecast Comparison" sheet: insert a new "Week_Start_Date" column
#     after "Week", renumber the Week labels (W01 -> W1, etc.), refresh the
#     MyForecast values, and store is_holiday_week as a boolean.
#   - "Summary" sheet: refresh the forecast totals that depend on the
#     corrected MyForecast numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date"); everything from the old
# column B onward shifts one column to the right (B->C, C->D, ... H->I,
# I->J).
$ws.Columns("B:B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Per-row data: Week label, Week start date, MyForecast, Amazon Mean,
# Amazon P70, Amazon P80, Amazon P90, is_holiday_week.
$rows = @(
  @{ Row=2;  Week="W1";  Date="2025-01-05"; My=110; Mean=168; P70=202; P80=238; P90=295; Holiday=$false },
  @{ Row=3;  Week="W2";  Date="2025-01-12"; My=96;  Mean=135; P70=165; P80=203; P90=264; Holiday=$false },
  @{ Row=4;  Week="W3";  Date="2025-01-19"; My=101; Mean=124; P70=152; P80=186; P90=242; Holiday=$false },
  @{ Row=5;  Week="W4";  Date="2025-01-26"; My=102; Mean=114; P70=139; P80=167; P90=213; Holiday=$false },
  @{ Row=6;  Week="W5";  Date="2025-02-02"; My=84;  Mean=74;  P70=90;  P80=107; P90=134; Holiday=$false },
  @{ Row=7;  Week="W6";  Date="2025-02-09"; My=77;  Mean=72;  P70=87;  P80=104; P90=130; Holiday=$false },
  @{ Row=8;  Week="W7";  Date="2025-02-16"; My=80;  Mean=69;  P70=84;  P80=100; P90=127; Holiday=$false },
  @{ Row=9;  Week="W8";  Date="2025-02-23"; My=69;  Mean=70;  P70=86;  P80=104; P90=133; Holiday=$false },
  @{ Row=10; Week="W9";  Date="2025-03-02"; My=69;  Mean=70;  P70=84;  P80=101; P90=126; Holiday=$false },
  @{ Row=11; Week="W10"; Date="2025-03-09"; My=67;  Mean=67;  P70=81;  P80=98;  P90=125; Holiday=$false },
  @{ Row=12; Week="W11"; Date="2025-03-16"; My=62;  Mean=61;  P70=75;  P80=92;  P90=121; Holiday=$false },
  @{ Row=13; Week="W12"; Date="2025-03-23"; My=76;  Mean=65;  P70=79;  P80=97;  P90=126; Holiday=$false },
  @{ Row=14; Week="W13"; Date="2025-03-30"; My=68;  Mean=65;  P70=79;  P80=95;  P90=121; Holiday=$false },
  @{ Row=15; Week="W14"; Date="2025-04-06"; My=70;  Mean=59;  P70=73;  P80=89;  P90=117; Holiday=$false },
  @{ Row=16; Week="W15"; Date="2025-04-13"; My=59;  Mean=58;  P70=72;  P80=88;  P90=116; Holiday=$false },
  @{ Row=17; Week="W16"; Date="2025-04-20"; My=58;  Mean=57;  P70=69;  P80=86;  P90=112; Holiday=$false }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Week label without the leading zero (W01 -> W1, ... W09 -> W9; W10+
    # already had no leading zero).
    $ws.Cells.Item($row, 1).Value = $r.Week

    # Week_Start_Date: must stay plain text, not be auto-converted to a
    # serial date by Excel's type inference.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.Style = "Normal"

    # ASIN (column C) is unchanged, already shifted by the column insert.

    # MyForecast (refreshed numbers).
    $ws.Cells.Item($row, 4).Value = $r.My
    $ws.Cells.Item($row, 5).Value = $r.Mean
    $ws.Cells.Item($row, 6).Value = $r.P70
    $ws.Cells.Item($row, 7).Value = $r.P80
    $ws.Cells.Item($row, 8).Value = $r.P90

    # Product Title (column I) is unchanged, already shifted by the
    # column insert.

    # is_holiday_week is now a boolean instead of a number.
    $ws.Cells.Item($row, 10).Value = $r.Holiday
}

# ---------------------------------------------------------------------
# Sheet 2: "Summary" -- refresh forecast totals
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

# These totals are stored as plain text (same as the rest of the Value
# column on this sheet), so they must not be auto-converted to numbers
# by Excel's type inference.
$totals = @(
  @{ Row=9;  Value="1249" },
  @{ Row=10; Value="719" },
  @{ Row=11; Value="409" }
)
foreach ($t in $totals) {
    $cell = $summary.Cells.Item($t.Row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $t.Value
    $cell.Style = "Normal"
}
